$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Corrected data values (column AI) - column AH holds formulas
#    (=AI+AE) that recalc automatically once AI changes.
# -----------------------------------------------------------------
$ws.Range("AI2").Value  = 460
$ws.Range("AI3").Value  = 460
$ws.Range("AI4").Value  = 460
$ws.Range("AI5").Value  = 400
$ws.Range("AI6").Value  = 400
$ws.Range("AI7").Value  = 414.5
$ws.Range("AI8").Value  = 414.5
$ws.Range("AI15").Value = 1289
$ws.Range("AI16").Value = 1289

# AI16 used to carry an explicit (no-op) border override; clear it so
# the cell falls back to the plain fill-only style used elsewhere in
# the column.
$ws.Range("AI16").Borders.LineStyle = -4142   # xlLineStyleNone

# -----------------------------------------------------------------
# 2) Number-format correction: the "0.00" (numFmtId 2) cells in this
#    block should use the workbook's "0.0" custom format (numFmtId
#    164) instead - apply to AJ2 and the AO (total) column.
# -----------------------------------------------------------------
$fmt = "0.0"

$ws.Range("AJ2").NumberFormat  = $fmt
$ws.Range("AO2").NumberFormat  = $fmt
$ws.Range("AO3").NumberFormat  = $fmt
$ws.Range("AO4").NumberFormat  = $fmt
$ws.Range("AO15").NumberFormat = $fmt
$ws.Range("AO16").NumberFormat = $fmt

$ws.Range("AO8").NumberFormat  = $fmt
$ws.Range("AO10").NumberFormat = $fmt
$ws.Range("AO11").NumberFormat = $fmt
$ws.Range("AO12").NumberFormat = $fmt
$ws.Range("AO13").NumberFormat = $fmt

$ws.Range("AO5").NumberFormat  = $fmt
$ws.Range("AO7").NumberFormat  = $fmt
$ws.Range("AO9").NumberFormat  = $fmt
$ws.Range("AO14").NumberFormat = $fmt

$ws.Range("AO6").NumberFormat  = $fmt

# -----------------------------------------------------------------
# 3) View state: scroll the sheet over and move the active selection.
# -----------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 35
$win.ScrollRow = 1
$ws.Range("AO7").Select()
